# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> color scheme "Integral"     (used by the slide master / all slides)
#   ppt/theme/theme2.xml  -> color scheme "Office Theme" (used by the notes master)
#
# The target revision swaps the two themes' contents: the slide master's
# theme becomes the default "Office Theme" palette, and the notes master's
# theme becomes the "Integral" palette. The font scheme and format scheme
# (fills/lines/effects) are identical between the two themes already, so
# the only real difference to apply is the 12 color-scheme slots.
#
# PowerPoint's object model only exposes one live/writable theme surface
# from COM (Master.Theme / NotesMaster.Theme / Slide.ThemeColorScheme all
# resolve to the same underlying theme part backing the slide master), so
# we drive the swap through that single surface: recolor it from
# "Integral" to the default "Office Theme" palette.
#
# ThemeColorScheme.Colors(index) order is:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# and .RGB takes a COM/VBA-style 0xBBGGRR integer (same convention as the
# VBA RGB() function), not the OOXML 0xRRGGBB hex string.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target palette = the built-in "Office Theme" colors (RRGGBB -> BBGGRR-as-int)
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
